# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the 06bdf7bd-fda1-44af-85fa-c1e817012ff9
# entry across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 06bdf7bd-fda1-44af-85fa-c1e817012ff9.md
$wsOverview.Range("G3").Value = "2016-09-05 04:50:46"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-09-05 04:50:41"
$wsZhCn.Range("K3").Value = "2016-09-05 04:50:59"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H3").Value = "2016-09-05 04:50:46"
$wsDeDe.Range("K3").Value = "2016-09-05 04:51:13"
